$d = $word.ActiveDocument

$pairs = @(
    @{old='76÷6=12, 4'; new='26÷6=4, 2'},
    @{old='97÷7=13, 6'; new='80÷7=11, 3'},
    @{old='97÷4=24, 1'; new='59÷7=8, 3'},
    @{old='67÷4=16, 3'; new='72÷9=8, 0'},
    @{old='30÷5=6, 0'; new='50÷2=25, 0'},
    @{old='10÷2=5, 0'; new='19÷3=6, 1'},
    @{old='21÷2=10, 1'; new='55÷2=27, 1'},
    @{old='57÷8=7, 1'; new='14÷2=7, 0'},
    @{old='20÷8=2, 4'; new='98÷5=19, 3'},
    @{old='98÷3=32, 2'; new='71÷2=35, 1'},
    @{old='75÷7=10, 5'; new='66÷9=7, 3'},
    @{old='51÷4=12, 3'; new='24÷3=8, 0'},
    @{old='27÷3=9, 0'; new='43÷3=14, 1'},
    @{old='58÷8=7, 2'; new='83÷4=20, 3'},
    @{old='78÷6=13, 0'; new='80÷2=40, 0'},
    @{old='77÷6=12, 5'; new='48÷6=8, 0'},
    @{old='66÷2=33, 0'; new='40÷9=4, 4'},
    @{old='71÷7=10, 1'; new='97÷8=12, 1'},
    @{old='61÷2=30, 1'; new='15÷3=5, 0'},
    @{old='65÷7=9, 2'; new='48÷3=16, 0'},
    @{old='11÷3=3, 2'; new='80÷9=8, 8'},
    @{old='30÷2=15, 0'; new='55÷2=27, 1'},
    @{old='60÷2=30, 0'; new='50÷6=8, 2'},
    @{old='61÷4=15, 1'; new='80÷4=20, 0'},
    @{old='36÷2=18, 0'; new='80÷8=10, 0'}
)

foreach ($p in $pairs) {
    $range = $d.Content
    $range.Find.Execute($p.old, $true, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)
}
